$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "27.492.78"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "1.872.60"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.73%  "
Set-TextValue $ws.Range("D5") "313.53"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("E6").Value = "  +0.76%  "
Set-TextValue $ws.Range("D7") "0.4799"
Set-TextValue $ws.Range("D8") "0.3772"
$ws.Range("E8").Value = "  +3.10%  "
Set-TextValue $ws.Range("D9") "0.07391"
$ws.Range("E9").Value = "  +2.92%  "
Set-TextValue $ws.Range("D10") "0.9403"
$ws.Range("E10").Value = "  +1.92%  "
Set-TextValue $ws.Range("D11") "20.69"
$ws.Range("E11").Value = "  +5.58%  "
Set-TextValue $ws.Range("D12") "0.07885"
$ws.Range("E12").Value = "  +3.76%  "
$ws.Range("D13").Value = "1.884.65"
$ws.Range("E13").Value = "  +2.95%  "
Set-TextValue $ws.Range("D14") "5.448"
$ws.Range("E14").Value = "  +2.81%  "
Set-TextValue $ws.Range("D15") "6.609"
$ws.Range("E15").Value = "  +3.38%  "
Set-TextValue $ws.Range("D16") "91.05"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("E17").Value = "  +0.66%  "
Set-TextValue $ws.Range("D18") "0.000008954"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("E19").Value = "  +0.69%  "
Set-TextValue $ws.Range("D20") "14.93"
$ws.Range("E20").Value = "  +2.75%  "
$ws.Range("D21").Value = "27.532.10"
$ws.Range("E21").Value = "  +2.25%  "
Set-TextValue $ws.Range("D22") "5.143"
$ws.Range("E22").Value = "  +2.41%  "
Set-TextValue $ws.Range("D24") "1.963"
$ws.Range("E24").Value = "  +2.48%  "
Set-TextValue $ws.Range("D25") "154.02"
$ws.Range("E25").Value = "  +1.15%  "
Set-TextValue $ws.Range("D26") "18.59"
$ws.Range("E26").Value = "  +2.32%  "
Set-TextValue $ws.Range("D27") "2.018"
$ws.Range("E27").Value = "  +0.72%  "
Set-TextValue $ws.Range("D28") "116.09"
$ws.Range("E28").Value = "  +1.53%  "
Set-TextValue $ws.Range("D29") "5.012"
$ws.Range("E29").Value = "  +2.96%  "
Set-TextValue $ws.Range("D30") "0.08931"
$ws.Range("E30").Value = "  +0.88%  "
Set-TextValue $ws.Range("D31") "3.329"
$ws.Range("E31").Value = "  +1.64%  "
Set-TextValue $ws.Range("D32") "1.219"
$ws.Range("E32").Value = "  +4.78%  "
Set-TextValue $ws.Range("D33") "4.604"
$ws.Range("E33").Value = "  +2.55%  "
Set-TextValue $ws.Range("D34") "0.7492"
$ws.Range("E34").Value = "  +0.45%  "
Set-TextValue $ws.Range("D35") "2.699"
$ws.Range("E35").Value = "  -3.32%  "
Set-TextValue $ws.Range("D36") "0.02073"
$ws.Range("E36").Value = "  +6.28%  "
Set-TextValue $ws.Range("D37") "1.122"
$ws.Range("E37").Value = "  +2.94%  "
Set-TextValue $ws.Range("D38") "0.05300"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("E39").Value = "  +1.11%  "
Set-TextValue $ws.Range("D40") "0.5383"
$ws.Range("E40").Value = "  +3.25%  "
Set-TextValue $ws.Range("D41") "7.096"
$ws.Range("E41").Value = "  +2.93%  "
Set-TextValue $ws.Range("D42") "0.1529"
$ws.Range("E42").Value = "  +1.20%  "
Set-TextValue $ws.Range("D43") "8.438"
$ws.Range("E43").Value = "  +3.16%  "
Set-TextValue $ws.Range("D44") "0.4851"
$ws.Range("E44").Value = "  +3.39%  "
Set-TextValue $ws.Range("D45") "10.63"
$ws.Range("E45").Value = "  +1.47%  "
Set-TextValue $ws.Range("D46") "1.016"
$ws.Range("E46").Value = "  +0.80%  "
Set-TextValue $ws.Range("D47") "1.667"
$ws.Range("E47").Value = "  +4.08%  "
Set-TextValue $ws.Range("D48") "103.11"
$ws.Range("E48").Value = "  +1.10%  "
Set-TextValue $ws.Range("D49") "67.22"
$ws.Range("E49").Value = "  +2.70%  "
Set-TextValue $ws.Range("D50") "0.06114"
Set-TextValue $ws.Range("D51") "0.9013"
$ws.Range("E51").Value = "  +1.81%  "
